$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.755.35'
$ws.Range('E2').Value = '  +0.45%  '
$ws.Range('D3').Value = '1.639.98'
$ws.Range('E3').Value = '  -0.13%  '
$ws.Range('E4').Value = '  +0.23%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '217.72'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.18%  '
$ws.Range('E7').Value = '  +0.30%  '
$ws.Range('E8').Value = '  -0.10%  '
$ws.Range('E9').Value = '  -0.28%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.09'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.07%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0844'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.23%  '
$ws.Range('D12').Value = '1.868.18'
$ws.Range('E12').Value = '  -0.14%  '
$ws.Range('D13').Value = '1.643.31'
$ws.Range('E13').Value = '  +0.16%  '
$ws.Range('E14').Value = '  -0.43%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.59'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.43%  '
$ws.Range('D17').Value = '26.741.68'
$ws.Range('E17').Value = '  +0.37%  '
$ws.Range('E18').Value = '  -1.80%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '212.71'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.51%  '
$ws.Range('E20').Value = '  +0.33%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.36'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.29%  '
$ws.Range('E22').Value = '  -0.83%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.33'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +4.20%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.25'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.30%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '145.59'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.14%  '
$ws.Range('E26').Value = '  +0.03%  '
$ws.Range('E27').Value = '  -1.53%  '
$ws.Range('E28').Value = '  +0.34%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.63'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.37%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0506'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.10%  '
$ws.Range('E32').Value = '  +1.15%  '
$ws.Range('E33').Value = '  +0.05%  '
$ws.Range('D34').Value = '1.278.89'
$ws.Range('E34').Value = '  +0.28%  '
$ws.Range('E35').Value = '  -0.41%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.44'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.71%  '
$ws.Range('E37').Value = '  -1.18%  '
$ws.Range('E38').Value = '  -0.05%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.814'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.53%  '
$ws.Range('E41').Value = '  -0.63%  '
$ws.Range('E42').Value = '  -1.56%  '
$ws.Range('D43').Value = '1.777.44'
$ws.Range('E43').Value = '  -0.18%  '
$ws.Range('E44').Value = '  -3.11%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '60.91'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +3.06%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '91.20'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.37%  '
$ws.Range('E47').Value = '  -1.02%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0520'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.02%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.57'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.25%  '
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0962'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.23%  '
$ws.Range('B51').Value = 'Mantle'
$ws.Range('C51').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.407'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.12%  '
